$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order shuffled match rows (columns F:V) ---
# Row 10
$ws.Cells.Item(10,"F").Value = 'Stolem Gniewino'
$ws.Cells.Item(10,"G").Value = 0
$ws.Cells.Item(10,"H").Value = 'Swinoujscie'
$ws.Cells.Item(10,"I").Value = 0
$ws.Cells.Item(10,"J").Value = 2.09
$ws.Cells.Item(10,"K").Value = '11/08/2023 17:45'
$ws.Cells.Item(10,"L").Value = 2.3
$ws.Cells.Item(10,"M").Value = '11/08/2023 17:59'
$ws.Cells.Item(10,"N").Value = 3.55
$ws.Cells.Item(10,"O").Value = '11/08/2023 17:45'
$ws.Cells.Item(10,"P").Value = 3.49
$ws.Cells.Item(10,"Q").Value = '11/08/2023 17:47'
$ws.Cells.Item(10,"R").Value = 2.91
$ws.Cells.Item(10,"S").Value = '11/08/2023 17:45'
$ws.Cells.Item(10,"T").Value = 2.61
$ws.Cells.Item(10,"U").Value = '11/08/2023 17:59'
$ws.Cells.Item(10,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/stolem-gniewino-swinoujscie/bqdM2Q3s/'

# Row 11
$ws.Cells.Item(11,"F").Value = 'Sroda'
$ws.Cells.Item(11,"G").Value = 1
$ws.Cells.Item(11,"H").Value = 'Notec Czarnkow'
$ws.Cells.Item(11,"I").Value = 4
$ws.Cells.Item(11,"J").Value = 1.58
$ws.Cells.Item(11,"K").Value = '11/08/2023 17:45'
$ws.Cells.Item(11,"L").Value = 1.41
$ws.Cells.Item(11,"M").Value = '11/08/2023 17:48'
$ws.Cells.Item(11,"N").Value = 4.12
$ws.Cells.Item(11,"O").Value = '11/08/2023 17:45'
$ws.Cells.Item(11,"P").Value = 4.57
$ws.Cells.Item(11,"Q").Value = '11/08/2023 17:48'
$ws.Cells.Item(11,"R").Value = 4.33
$ws.Cells.Item(11,"S").Value = '11/08/2023 17:45'
$ws.Cells.Item(11,"T").Value = 5.56
$ws.Cells.Item(11,"U").Value = '11/08/2023 17:48'
$ws.Cells.Item(11,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/polonia-sroda-wielkopol-notec-czarnkow/SOmwar36/'

# Row 12
$ws.Cells.Item(12,"F").Value = 'Swit Skolwin'
$ws.Cells.Item(12,"G").Value = 2
$ws.Cells.Item(12,"H").Value = 'Gedania Gdansk'
$ws.Cells.Item(12,"I").Value = 0
$ws.Cells.Item(12,"J").Value = 1.74
$ws.Cells.Item(12,"K").Value = '10/08/2023 06:12'
$ws.Cells.Item(12,"L").Value = 1.67
$ws.Cells.Item(12,"M").Value = '11/08/2023 17:39'
$ws.Cells.Item(12,"N").Value = 3.6
$ws.Cells.Item(12,"O").Value = '10/08/2023 06:12'
$ws.Cells.Item(12,"P").Value = 3.86
$ws.Cells.Item(12,"Q").Value = '11/08/2023 17:39'
$ws.Cells.Item(12,"R").Value = 3.39
$ws.Cells.Item(12,"S").Value = '10/08/2023 06:12'
$ws.Cells.Item(12,"T").Value = 3.96
$ws.Cells.Item(12,"U").Value = '11/08/2023 17:39'
$ws.Cells.Item(12,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/swit-skolwin-gedania-gdansk/Momsb2IC/'

# Row 24
$ws.Cells.Item(24,"F").Value = 'Elana Torun'
$ws.Cells.Item(24,"G").Value = 1
$ws.Cells.Item(24,"H").Value = 'Vineta W.'
$ws.Cells.Item(24,"I").Value = 2
$ws.Cells.Item(24,"J").Value = 1.86
$ws.Cells.Item(24,"K").Value = '19/08/2023 09:25'
$ws.Cells.Item(24,"L").Value = 1.83
$ws.Cells.Item(24,"M").Value = '19/08/2023 14:32'
$ws.Cells.Item(24,"N").Value = 3.54
$ws.Cells.Item(24,"O").Value = '19/08/2023 09:25'
$ws.Cells.Item(24,"P").Value = 3.49
$ws.Cells.Item(24,"Q").Value = '19/08/2023 15:45'
$ws.Cells.Item(24,"R").Value = 3.41
$ws.Cells.Item(24,"S").Value = '19/08/2023 09:25'
$ws.Cells.Item(24,"T").Value = 3.67
$ws.Cells.Item(24,"U").Value = '19/08/2023 15:45'
$ws.Cells.Item(24,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/elana-torun-vineta-wolin/69LaGqYt/'

# Row 25
$ws.Cells.Item(25,"F").Value = 'Swinoujscie'
$ws.Cells.Item(25,"G").Value = 0
$ws.Cells.Item(25,"H").Value = 'Blekitni Stargard'
$ws.Cells.Item(25,"I").Value = 3
$ws.Cells.Item(25,"J").Value = 2.32
$ws.Cells.Item(25,"K").Value = '19/08/2023 09:26'
$ws.Cells.Item(25,"L").Value = 2.34
$ws.Cells.Item(25,"M").Value = '19/08/2023 16:59'
$ws.Cells.Item(25,"N").Value = 3.3
$ws.Cells.Item(25,"O").Value = '19/08/2023 09:26'
$ws.Cells.Item(25,"P").Value = 3.33
$ws.Cells.Item(25,"Q").Value = '19/08/2023 16:59'
$ws.Cells.Item(25,"R").Value = 2.64
$ws.Cells.Item(25,"S").Value = '19/08/2023 09:26'
$ws.Cells.Item(25,"T").Value = 2.66
$ws.Cells.Item(25,"U").Value = '19/08/2023 16:59'
$ws.Cells.Item(25,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/swinoujscie-blekitni-stargard/SO98EN2g/'

# Row 31
$ws.Cells.Item(31,"F").Value = 'Stolem Gniewino'
$ws.Cells.Item(31,"G").Value = 0
$ws.Cells.Item(31,"H").Value = 'Cartusia Kartuzy'
$ws.Cells.Item(31,"I").Value = 0
$ws.Cells.Item(31,"J").Value = 2.62
$ws.Cells.Item(31,"K").Value = '22/08/2023 06:12'
$ws.Cells.Item(31,"L").Value = 2.64
$ws.Cells.Item(31,"M").Value = '23/08/2023 17:59'
$ws.Cells.Item(31,"N").Value = 3.11
$ws.Cells.Item(31,"O").Value = '22/08/2023 06:12'
$ws.Cells.Item(31,"P").Value = 3.56
$ws.Cells.Item(31,"Q").Value = '23/08/2023 17:59'
$ws.Cells.Item(31,"R").Value = 2.26
$ws.Cells.Item(31,"S").Value = '22/08/2023 06:12'
$ws.Cells.Item(31,"T").Value = 2.25
$ws.Cells.Item(31,"U").Value = '23/08/2023 17:59'
$ws.Cells.Item(31,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/stolem-gniewino-cartusia-kartuzy/noG3YrfP/'

# Row 32
$ws.Cells.Item(32,"F").Value = 'Vineta W.'
$ws.Cells.Item(32,"G").Value = 2
$ws.Cells.Item(32,"H").Value = 'Solec Kujawski'
$ws.Cells.Item(32,"I").Value = 1
$ws.Cells.Item(32,"J").Value = 1.76
$ws.Cells.Item(32,"K").Value = '22/08/2023 06:12'
$ws.Cells.Item(32,"L").Value = 1.74
$ws.Cells.Item(32,"M").Value = '23/08/2023 17:00'
$ws.Cells.Item(32,"N").Value = 3.54
$ws.Cells.Item(32,"O").Value = '22/08/2023 06:12'
$ws.Cells.Item(32,"P").Value = 3.85
$ws.Cells.Item(32,"Q").Value = '23/08/2023 17:00'
$ws.Cells.Item(32,"R").Value = 3.34
$ws.Cells.Item(32,"S").Value = '22/08/2023 06:12'
$ws.Cells.Item(32,"T").Value = 3.68
$ws.Cells.Item(32,"U").Value = '23/08/2023 17:00'
$ws.Cells.Item(32,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/vineta-wolin-unia-solec-kujawski/pMevStXh/'

# Row 33
$ws.Cells.Item(33,"F").Value = 'Zawisza'
$ws.Cells.Item(33,"G").Value = 2
$ws.Cells.Item(33,"H").Value = 'Swinoujscie'
$ws.Cells.Item(33,"I").Value = 0
$ws.Cells.Item(33,"J").Value = 1.21
$ws.Cells.Item(33,"K").Value = '23/08/2023 11:12'
$ws.Cells.Item(33,"L").Value = 1.31
$ws.Cells.Item(33,"M").Value = '23/08/2023 17:02'
$ws.Cells.Item(33,"N").Value = 5.94
$ws.Cells.Item(33,"O").Value = '23/08/2023 11:12'
$ws.Cells.Item(33,"P").Value = 6.41
$ws.Cells.Item(33,"Q").Value = '23/08/2023 17:02'
$ws.Cells.Item(33,"R").Value = 8.15
$ws.Cells.Item(33,"S").Value = '23/08/2023 11:12'
$ws.Cells.Item(33,"T").Value = 5.43
$ws.Cells.Item(33,"U").Value = '23/08/2023 17:02'
$ws.Cells.Item(33,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/zawisza-swinoujscie/6VdzTMHn/'

# Row 34
$ws.Cells.Item(34,"F").Value = 'Swit Skolwin'
$ws.Cells.Item(34,"G").Value = 3
$ws.Cells.Item(34,"H").Value = 'Luzino'
$ws.Cells.Item(34,"I").Value = 0
$ws.Cells.Item(34,"J").Value = 1.2
$ws.Cells.Item(34,"K").Value = '23/08/2023 11:12'
$ws.Cells.Item(34,"L").Value = 1.28
$ws.Cells.Item(34,"M").Value = '23/08/2023 17:50'
$ws.Cells.Item(34,"N").Value = 6.33
$ws.Cells.Item(34,"O").Value = '23/08/2023 11:12'
$ws.Cells.Item(34,"P").Value = 5.76
$ws.Cells.Item(34,"Q").Value = '23/08/2023 17:50'
$ws.Cells.Item(34,"R").Value = 7.73
$ws.Cells.Item(34,"S").Value = '23/08/2023 11:12'
$ws.Cells.Item(34,"T").Value = 6.55
$ws.Cells.Item(34,"U").Value = '23/08/2023 17:50'
$ws.Cells.Item(34,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/swit-skolwin-luzino/hnanQK15/'

# Row 35
$ws.Cells.Item(35,"F").Value = 'Nowe Skalmierzyce'
$ws.Cells.Item(35,"G").Value = 0
$ws.Cells.Item(35,"H").Value = 'Pogon Szczecin II'
$ws.Cells.Item(35,"I").Value = 5
$ws.Cells.Item(35,"J").Value = 2.6
$ws.Cells.Item(35,"K").Value = '22/08/2023 06:12'
$ws.Cells.Item(35,"L").Value = 1.95
$ws.Cells.Item(35,"M").Value = '23/08/2023 17:45'
$ws.Cells.Item(35,"N").Value = 3.2
$ws.Cells.Item(35,"O").Value = '22/08/2023 06:12'
$ws.Cells.Item(35,"P").Value = 3.4
$ws.Cells.Item(35,"Q").Value = '23/08/2023 17:45'
$ws.Cells.Item(35,"R").Value = 2.22
$ws.Cells.Item(35,"S").Value = '22/08/2023 06:12'
$ws.Cells.Item(35,"T").Value = 2.9
$ws.Cells.Item(35,"U").Value = '23/08/2023 17:45'
$ws.Cells.Item(35,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/nowe-skalmierzyce-pogon-szczecin/YsSL9aHO/'

# Row 58
$ws.Cells.Item(58,"F").Value = 'Swinoujscie'
$ws.Cells.Item(58,"G").Value = 1
$ws.Cells.Item(58,"H").Value = 'Swit Skolwin'
$ws.Cells.Item(58,"I").Value = 3
$ws.Cells.Item(58,"J").Value = 3.32
$ws.Cells.Item(58,"K").Value = '08/09/2023 05:13'
$ws.Cells.Item(58,"L").Value = 3.75
$ws.Cells.Item(58,"M").Value = '09/09/2023 16:57'
$ws.Cells.Item(58,"N").Value = 3.44
$ws.Cells.Item(58,"O").Value = '08/09/2023 05:13'
$ws.Cells.Item(58,"P").Value = 3.5
$ws.Cells.Item(58,"Q").Value = '09/09/2023 16:57'
$ws.Cells.Item(58,"R").Value = 1.83
$ws.Cells.Item(58,"S").Value = '08/09/2023 05:13'
$ws.Cells.Item(58,"T").Value = 1.81
$ws.Cells.Item(58,"U").Value = '09/09/2023 16:57'
$ws.Cells.Item(58,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/swinoujscie-swit-skolwin/2c8mWnUe/'

# Row 59
$ws.Cells.Item(59,"F").Value = 'Starogard Gdanski'
$ws.Cells.Item(59,"G").Value = 1
$ws.Cells.Item(59,"H").Value = 'Blekitni Stargard'
$ws.Cells.Item(59,"I").Value = 2
$ws.Cells.Item(59,"J").Value = 1.89
$ws.Cells.Item(59,"K").Value = '08/09/2023 05:13'
$ws.Cells.Item(59,"L").Value = 2
$ws.Cells.Item(59,"M").Value = '09/09/2023 16:23'
$ws.Cells.Item(59,"N").Value = 3.41
$ws.Cells.Item(59,"O").Value = '08/09/2023 05:13'
$ws.Cells.Item(59,"P").Value = 3.78
$ws.Cells.Item(59,"Q").Value = '09/09/2023 16:22'
$ws.Cells.Item(59,"R").Value = 3.05
$ws.Cells.Item(59,"S").Value = '08/09/2023 05:13'
$ws.Cells.Item(59,"T").Value = 2.94
$ws.Cells.Item(59,"U").Value = '09/09/2023 16:23'
$ws.Cells.Item(59,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/starogard-gdanski-blekitni-stargard/jJlrhUjF/'

# Row 74
$ws.Cells.Item(74,"F").Value = 'Stolem Gniewino'
$ws.Cells.Item(74,"G").Value = 0
$ws.Cells.Item(74,"H").Value = 'Blekitni Stargard'
$ws.Cells.Item(74,"I").Value = 6
$ws.Cells.Item(74,"J").Value = 2.23
$ws.Cells.Item(74,"K").Value = '22/09/2023 03:13'
$ws.Cells.Item(74,"L").Value = 2.53
$ws.Cells.Item(74,"M").Value = '23/09/2023 15:42'
$ws.Cells.Item(74,"N").Value = 3.26
$ws.Cells.Item(74,"O").Value = '22/09/2023 03:13'
$ws.Cells.Item(74,"P").Value = 3.34
$ws.Cells.Item(74,"Q").Value = '23/09/2023 15:42'
$ws.Cells.Item(74,"R").Value = 2.56
$ws.Cells.Item(74,"S").Value = '22/09/2023 03:13'
$ws.Cells.Item(74,"T").Value = 2.44
$ws.Cells.Item(74,"U").Value = '23/09/2023 15:42'
$ws.Cells.Item(74,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/stolem-gniewino-blekitni-stargard/Uq4eIQ4K/'

# Row 75
$ws.Cells.Item(75,"F").Value = 'Swinoujscie'
$ws.Cells.Item(75,"G").Value = 1
$ws.Cells.Item(75,"H").Value = 'Notec Czarnkow'
$ws.Cells.Item(75,"I").Value = 3
$ws.Cells.Item(75,"J").Value = 2.1
$ws.Cells.Item(75,"K").Value = '22/09/2023 03:13'
$ws.Cells.Item(75,"L").Value = 2.3
$ws.Cells.Item(75,"M").Value = '23/09/2023 15:43'
$ws.Cells.Item(75,"N").Value = 3.36
$ws.Cells.Item(75,"O").Value = '22/09/2023 03:13'
$ws.Cells.Item(75,"P").Value = 3.88
$ws.Cells.Item(75,"Q").Value = '23/09/2023 15:43'
$ws.Cells.Item(75,"R").Value = 2.69
$ws.Cells.Item(75,"S").Value = '22/09/2023 03:13'
$ws.Cells.Item(75,"T").Value = 2.39
$ws.Cells.Item(75,"U").Value = '23/09/2023 15:43'
$ws.Cells.Item(75,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/swinoujscie-notec-czarnkow/pdUib5Dl/'

# Row 76
$ws.Cells.Item(76,"F").Value = 'Elana Torun'
$ws.Cells.Item(76,"G").Value = 1
$ws.Cells.Item(76,"H").Value = 'Kleczew'
$ws.Cells.Item(76,"I").Value = 0
$ws.Cells.Item(76,"J").Value = 1.75
$ws.Cells.Item(76,"K").Value = '22/09/2023 03:13'
$ws.Cells.Item(76,"L").Value = 1.84
$ws.Cells.Item(76,"M").Value = '23/09/2023 15:52'
$ws.Cells.Item(76,"N").Value = 3.58
$ws.Cells.Item(76,"O").Value = '22/09/2023 03:13'
$ws.Cells.Item(76,"P").Value = 3.41
$ws.Cells.Item(76,"Q").Value = '23/09/2023 15:53'
$ws.Cells.Item(76,"R").Value = 3.35
$ws.Cells.Item(76,"S").Value = '22/09/2023 03:13'
$ws.Cells.Item(76,"T").Value = 3.73
$ws.Cells.Item(76,"U").Value = '23/09/2023 15:52'
$ws.Cells.Item(76,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/elana-torun-kleczew/lQrM3RkR/'

# Row 77
$ws.Cells.Item(77,"F").Value = 'Starogard Gdanski'
$ws.Cells.Item(77,"G").Value = 2
$ws.Cells.Item(77,"H").Value = 'Vineta W.'
$ws.Cells.Item(77,"I").Value = 3
$ws.Cells.Item(77,"J").Value = 2.37
$ws.Cells.Item(77,"K").Value = '22/09/2023 05:53'
$ws.Cells.Item(77,"L").Value = 3.36
$ws.Cells.Item(77,"M").Value = '23/09/2023 15:53'
$ws.Cells.Item(77,"N").Value = 3.45
$ws.Cells.Item(77,"O").Value = '22/09/2023 05:53'
$ws.Cells.Item(77,"P").Value = 3.69
$ws.Cells.Item(77,"Q").Value = '23/09/2023 15:53'
$ws.Cells.Item(77,"R").Value = 2.39
$ws.Cells.Item(77,"S").Value = '22/09/2023 05:53'
$ws.Cells.Item(77,"T").Value = 1.86
$ws.Cells.Item(77,"U").Value = '23/09/2023 15:53'
$ws.Cells.Item(77,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/starogard-gdanski-vineta-wolin/llR8fNCD/'

# Row 89
$ws.Cells.Item(89,"F").Value = 'Pogon Szczecin II'
$ws.Cells.Item(89,"G").Value = 5
$ws.Cells.Item(89,"H").Value = 'Gedania Gdansk'
$ws.Cells.Item(89,"I").Value = 2
$ws.Cells.Item(89,"J").Value = 1.54
$ws.Cells.Item(89,"K").Value = '05/10/2023 23:12'
$ws.Cells.Item(89,"L").Value = 1.56
$ws.Cells.Item(89,"M").Value = '07/10/2023 11:40'
$ws.Cells.Item(89,"N").Value = 3.99
$ws.Cells.Item(89,"O").Value = '05/10/2023 23:12'
$ws.Cells.Item(89,"P").Value = 4.25
$ws.Cells.Item(89,"Q").Value = '07/10/2023 11:40'
$ws.Cells.Item(89,"R").Value = 3.98
$ws.Cells.Item(89,"S").Value = '05/10/2023 23:12'
$ws.Cells.Item(89,"T").Value = 4.32
$ws.Cells.Item(89,"U").Value = '07/10/2023 11:40'
$ws.Cells.Item(89,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/pogon-szczecin-gedania-gdansk/8CMGNLt0/'

# Row 90
$ws.Cells.Item(90,"F").Value = 'Solec Kujawski'
$ws.Cells.Item(90,"G").Value = 1
$ws.Cells.Item(90,"H").Value = 'Kleczew'
$ws.Cells.Item(90,"I").Value = 0
$ws.Cells.Item(90,"J").Value = 2.86
$ws.Cells.Item(90,"K").Value = '05/10/2023 23:12'
$ws.Cells.Item(90,"L").Value = 3.6
$ws.Cells.Item(90,"M").Value = '07/10/2023 11:54'
$ws.Cells.Item(90,"N").Value = 3.57
$ws.Cells.Item(90,"O").Value = '05/10/2023 23:12'
$ws.Cells.Item(90,"P").Value = 3.48
$ws.Cells.Item(90,"Q").Value = '07/10/2023 11:54'
$ws.Cells.Item(90,"R").Value = 1.93
$ws.Cells.Item(90,"S").Value = '05/10/2023 23:12'
$ws.Cells.Item(90,"T").Value = 1.85
$ws.Cells.Item(90,"U").Value = '07/10/2023 11:54'
$ws.Cells.Item(90,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/unia-solec-kujawski-kleczew/GEsxV3lQ/'

# Row 91
$ws.Cells.Item(91,"F").Value = 'Cartusia Kartuzy'
$ws.Cells.Item(91,"G").Value = 4
$ws.Cells.Item(91,"H").Value = 'Notec Czarnkow'
$ws.Cells.Item(91,"I").Value = 0
$ws.Cells.Item(91,"J").Value = 1.7
$ws.Cells.Item(91,"K").Value = '06/10/2023 00:12'
$ws.Cells.Item(91,"L").Value = 1.76
$ws.Cells.Item(91,"M").Value = '07/10/2023 12:57'
$ws.Cells.Item(91,"N").Value = 3.6
$ws.Cells.Item(91,"O").Value = '06/10/2023 00:12'
$ws.Cells.Item(91,"P").Value = 3.94
$ws.Cells.Item(91,"Q").Value = '07/10/2023 12:58'
$ws.Cells.Item(91,"R").Value = 3.54
$ws.Cells.Item(91,"S").Value = '06/10/2023 00:12'
$ws.Cells.Item(91,"T").Value = 3.52
$ws.Cells.Item(91,"U").Value = '07/10/2023 12:57'
$ws.Cells.Item(91,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/cartusia-kartuzy-notec-czarnkow/ULNCO1Rg/'

# Row 92
$ws.Cells.Item(92,"F").Value = 'Unia Swarzedz'
$ws.Cells.Item(92,"G").Value = 2
$ws.Cells.Item(92,"H").Value = 'Luzino'
$ws.Cells.Item(92,"I").Value = 1
$ws.Cells.Item(92,"J").Value = 1.53
$ws.Cells.Item(92,"K").Value = '06/10/2023 00:12'
$ws.Cells.Item(92,"L").Value = 1.57
$ws.Cells.Item(92,"M").Value = '07/10/2023 12:57'
$ws.Cells.Item(92,"N").Value = 3.96
$ws.Cells.Item(92,"O").Value = '06/10/2023 00:12'
$ws.Cells.Item(92,"P").Value = 4.38
$ws.Cells.Item(92,"Q").Value = '07/10/2023 12:57'
$ws.Cells.Item(92,"R").Value = 4.09
$ws.Cells.Item(92,"S").Value = '06/10/2023 00:12'
$ws.Cells.Item(92,"T").Value = 4.11
$ws.Cells.Item(92,"U").Value = '07/10/2023 12:57'
$ws.Cells.Item(92,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/unia-swarzedz-luzino/0nJ8PsCm/'

# Row 94
$ws.Cells.Item(94,"F").Value = 'Stolem Gniewino'
$ws.Cells.Item(94,"G").Value = 2
$ws.Cells.Item(94,"H").Value = 'Vineta W.'
$ws.Cells.Item(94,"I").Value = 3
$ws.Cells.Item(94,"J").Value = 2.75
$ws.Cells.Item(94,"K").Value = '06/10/2023 02:12'
$ws.Cells.Item(94,"L").Value = 3.82
$ws.Cells.Item(94,"M").Value = '07/10/2023 14:05'
$ws.Cells.Item(94,"N").Value = 3.31
$ws.Cells.Item(94,"O").Value = '06/10/2023 02:12'
$ws.Cells.Item(94,"P").Value = 3.8
$ws.Cells.Item(94,"Q").Value = '07/10/2023 14:05'
$ws.Cells.Item(94,"R").Value = 2.07
$ws.Cells.Item(94,"S").Value = '06/10/2023 02:12'
$ws.Cells.Item(94,"T").Value = 1.72
$ws.Cells.Item(94,"U").Value = '07/10/2023 14:05'
$ws.Cells.Item(94,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/stolem-gniewino-vineta-wolin/b1FTKJQI/'

# Row 95
$ws.Cells.Item(95,"F").Value = 'Swinoujscie'
$ws.Cells.Item(95,"G").Value = 3
$ws.Cells.Item(95,"H").Value = 'Elana Torun'
$ws.Cells.Item(95,"I").Value = 0
$ws.Cells.Item(95,"J").Value = 3.06
$ws.Cells.Item(95,"K").Value = '06/10/2023 02:12'
$ws.Cells.Item(95,"L").Value = 4.52
$ws.Cells.Item(95,"M").Value = '07/10/2023 14:07'
$ws.Cells.Item(95,"N").Value = 3.42
$ws.Cells.Item(95,"O").Value = '06/10/2023 02:12'
$ws.Cells.Item(95,"P").Value = 3.77
$ws.Cells.Item(95,"Q").Value = '07/10/2023 14:07'
$ws.Cells.Item(95,"R").Value = 1.88
$ws.Cells.Item(95,"S").Value = '06/10/2023 02:12'
$ws.Cells.Item(95,"T").Value = 1.61
$ws.Cells.Item(95,"U").Value = '07/10/2023 14:07'
$ws.Cells.Item(95,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/swinoujscie-elana-torun/MeK4QNds/'

# Row 98
$ws.Cells.Item(98,"F").Value = 'Luzino'
$ws.Cells.Item(98,"G").Value = 3
$ws.Cells.Item(98,"H").Value = 'Cartusia Kartuzy'
$ws.Cells.Item(98,"I").Value = 1
$ws.Cells.Item(98,"J").Value = 4.48
$ws.Cells.Item(98,"K").Value = '13/10/2023 00:13'
$ws.Cells.Item(98,"L").Value = 5.93
$ws.Cells.Item(98,"M").Value = '14/10/2023 12:57'
$ws.Cells.Item(98,"N").Value = 3.92
$ws.Cells.Item(98,"O").Value = '13/10/2023 00:13'
$ws.Cells.Item(98,"P").Value = 4.42
$ws.Cells.Item(98,"Q").Value = '14/10/2023 12:57'
$ws.Cells.Item(98,"R").Value = 1.51
$ws.Cells.Item(98,"S").Value = '13/10/2023 00:13'
$ws.Cells.Item(98,"T").Value = 1.41
$ws.Cells.Item(98,"U").Value = '14/10/2023 12:57'
$ws.Cells.Item(98,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/luzino-cartusia-kartuzy/EiQf6uTM/'

# Row 99
$ws.Cells.Item(99,"F").Value = 'Unia Swarzedz'
$ws.Cells.Item(99,"G").Value = 0
$ws.Cells.Item(99,"H").Value = 'Elana Torun'
$ws.Cells.Item(99,"I").Value = 0
$ws.Cells.Item(99,"J").Value = 2.69
$ws.Cells.Item(99,"K").Value = '13/10/2023 00:13'
$ws.Cells.Item(99,"L").Value = 2.16
$ws.Cells.Item(99,"M").Value = '14/10/2023 12:56'
$ws.Cells.Item(99,"N").Value = 3.23
$ws.Cells.Item(99,"O").Value = '13/10/2023 00:13'
$ws.Cells.Item(99,"P").Value = 3.28
$ws.Cells.Item(99,"Q").Value = '14/10/2023 12:56'
$ws.Cells.Item(99,"R").Value = 2.14
$ws.Cells.Item(99,"S").Value = '13/10/2023 00:13'
$ws.Cells.Item(99,"T").Value = 2.97
$ws.Cells.Item(99,"U").Value = '14/10/2023 12:56'
$ws.Cells.Item(99,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/unia-swarzedz-elana-torun/hrPb5arT/'

# Row 108
$ws.Cells.Item(108,"F").Value = 'Nowe Skalmierzyce'
$ws.Cells.Item(108,"G").Value = 0
$ws.Cells.Item(108,"H").Value = 'Gedania Gdansk'
$ws.Cells.Item(108,"I").Value = 2
$ws.Cells.Item(108,"J").Value = 1.93
$ws.Cells.Item(108,"K").Value = '20/10/2023 01:13'
$ws.Cells.Item(108,"L").Value = 2.01
$ws.Cells.Item(108,"M").Value = '21/10/2023 13:42'
$ws.Cells.Item(108,"N").Value = 3.56
$ws.Cells.Item(108,"O").Value = '20/10/2023 01:13'
$ws.Cells.Item(108,"P").Value = 3.78
$ws.Cells.Item(108,"Q").Value = '21/10/2023 13:42'
$ws.Cells.Item(108,"R").Value = 2.85
$ws.Cells.Item(108,"S").Value = '20/10/2023 01:13'
$ws.Cells.Item(108,"T").Value = 2.92
$ws.Cells.Item(108,"U").Value = '21/10/2023 13:42'
$ws.Cells.Item(108,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/nowe-skalmierzyce-gedania-gdansk/hroUKvrG/'

# Row 109
$ws.Cells.Item(109,"F").Value = 'Zawisza'
$ws.Cells.Item(109,"G").Value = 2
$ws.Cells.Item(109,"H").Value = 'Vineta W.'
$ws.Cells.Item(109,"I").Value = 2
$ws.Cells.Item(109,"J").Value = 1.63
$ws.Cells.Item(109,"K").Value = '20/10/2023 01:13'
$ws.Cells.Item(109,"L").Value = 1.75
$ws.Cells.Item(109,"M").Value = '21/10/2023 13:10'
$ws.Cells.Item(109,"N").Value = 3.74
$ws.Cells.Item(109,"O").Value = '20/10/2023 01:13'
$ws.Cells.Item(109,"P").Value = 3.86
$ws.Cells.Item(109,"Q").Value = '21/10/2023 13:10'
$ws.Cells.Item(109,"R").Value = 3.71
$ws.Cells.Item(109,"S").Value = '20/10/2023 01:13'
$ws.Cells.Item(109,"T").Value = 3.63
$ws.Cells.Item(109,"U").Value = '21/10/2023 13:10'
$ws.Cells.Item(109,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/zawisza-vineta-wolin/jJQ4ExLq/'

# --- Append new rows 112-121 ---
$ws.Range("A111:V111").Copy()
$ws.Range("A112:V121").PasteSpecial(-4122)

# Row 112
$ws.Cells.Item(112,"A").Value = 111
$ws.Cells.Item(112,"B").Value = 'poland'
$ws.Cells.Item(112,"C").Value = 'iii-liga-group-ii'
$ws.Cells.Item(112,"D").Value = '2023-2024'
$ws.Cells.Item(112,"E").Value = 45227.5
$ws.Cells.Item(112,"F").Value = 'Gedania Gdansk'
$ws.Cells.Item(112,"G").Value = 0
$ws.Cells.Item(112,"H").Value = 'Stolem Gniewino'
$ws.Cells.Item(112,"I").Value = 1
$ws.Cells.Item(112,"J").Value = 1.43
$ws.Cells.Item(112,"K").Value = '27/10/2023 00:12'
$ws.Cells.Item(112,"L").Value = 1.7
$ws.Cells.Item(112,"M").Value = '28/10/2023 11:57'
$ws.Cells.Item(112,"N").Value = 4.34
$ws.Cells.Item(112,"O").Value = '27/10/2023 00:12'
$ws.Cells.Item(112,"P").Value = 3.95
$ws.Cells.Item(112,"Q").Value = '28/10/2023 11:57'
$ws.Cells.Item(112,"R").Value = 4.5
$ws.Cells.Item(112,"S").Value = '27/10/2023 00:12'
$ws.Cells.Item(112,"T").Value = 3.78
$ws.Cells.Item(112,"U").Value = '28/10/2023 11:57'
$ws.Cells.Item(112,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/gedania-gdansk-stolem-gniewino/vgMbzKs3/'

# Row 113
$ws.Cells.Item(113,"A").Value = 112
$ws.Cells.Item(113,"B").Value = 'poland'
$ws.Cells.Item(113,"C").Value = 'iii-liga-group-ii'
$ws.Cells.Item(113,"D").Value = '2023-2024'
$ws.Cells.Item(113,"E").Value = 45227.5
$ws.Cells.Item(113,"F").Value = 'Swit Skolwin'
$ws.Cells.Item(113,"G").Value = 3
$ws.Cells.Item(113,"H").Value = 'Blekitni Stargard'
$ws.Cells.Item(113,"I").Value = 0
$ws.Cells.Item(113,"J").Value = 1.44
$ws.Cells.Item(113,"K").Value = '27/10/2023 00:12'
$ws.Cells.Item(113,"L").Value = 1.5
$ws.Cells.Item(113,"M").Value = '28/10/2023 00:08'
$ws.Cells.Item(113,"N").Value = 4.06
$ws.Cells.Item(113,"O").Value = '27/10/2023 00:12'
$ws.Cells.Item(113,"P").Value = 4.09
$ws.Cells.Item(113,"Q").Value = '28/10/2023 10:01'
$ws.Cells.Item(113,"R").Value = 4.74
$ws.Cells.Item(113,"S").Value = '27/10/2023 00:12'
$ws.Cells.Item(113,"T").Value = 5.04
$ws.Cells.Item(113,"U").Value = '28/10/2023 06:04'
$ws.Cells.Item(113,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/swit-skolwin-blekitni-stargard/QkIfy0Sd/'

# Row 114
$ws.Cells.Item(114,"A").Value = 113
$ws.Cells.Item(114,"B").Value = 'poland'
$ws.Cells.Item(114,"C").Value = 'iii-liga-group-ii'
$ws.Cells.Item(114,"D").Value = '2023-2024'
$ws.Cells.Item(114,"E").Value = 45227.5
$ws.Cells.Item(114,"F").Value = 'Solec Kujawski'
$ws.Cells.Item(114,"G").Value = 1
$ws.Cells.Item(114,"H").Value = 'Cartusia Kartuzy'
$ws.Cells.Item(114,"I").Value = 6
$ws.Cells.Item(114,"J").Value = 3.6
$ws.Cells.Item(114,"K").Value = '27/10/2023 00:12'
$ws.Cells.Item(114,"L").Value = 4.32
$ws.Cells.Item(114,"M").Value = '28/10/2023 11:59'
$ws.Cells.Item(114,"N").Value = 3.64
$ws.Cells.Item(114,"O").Value = '27/10/2023 00:12'
$ws.Cells.Item(114,"P").Value = 3.8
$ws.Cells.Item(114,"Q").Value = '28/10/2023 11:59'
$ws.Cells.Item(114,"R").Value = 1.68
$ws.Cells.Item(114,"S").Value = '27/10/2023 00:12'
$ws.Cells.Item(114,"T").Value = 1.63
$ws.Cells.Item(114,"U").Value = '28/10/2023 11:59'
$ws.Cells.Item(114,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/unia-solec-kujawski-cartusia-kartuzy/8YDFWysS/'

# Row 115
$ws.Cells.Item(115,"A").Value = 114
$ws.Cells.Item(115,"B").Value = 'poland'
$ws.Cells.Item(115,"C").Value = 'iii-liga-group-ii'
$ws.Cells.Item(115,"D").Value = '2023-2024'
$ws.Cells.Item(115,"E").Value = 45227.54166666666
$ws.Cells.Item(115,"F").Value = 'Swinoujscie'
$ws.Cells.Item(115,"G").Value = 2
$ws.Cells.Item(115,"H").Value = 'Unia Swarzedz'
$ws.Cells.Item(115,"I").Value = 4
$ws.Cells.Item(115,"J").Value = 2.5
$ws.Cells.Item(115,"K").Value = '27/10/2023 01:12'
$ws.Cells.Item(115,"L").Value = 2.96
$ws.Cells.Item(115,"M").Value = '28/10/2023 12:59'
$ws.Cells.Item(115,"N").Value = 3.28
$ws.Cells.Item(115,"O").Value = '27/10/2023 01:12'
$ws.Cells.Item(115,"P").Value = 3.23
$ws.Cells.Item(115,"Q").Value = '28/10/2023 12:36'
$ws.Cells.Item(115,"R").Value = 2.26
$ws.Cells.Item(115,"S").Value = '27/10/2023 01:12'
$ws.Cells.Item(115,"T").Value = 2.19
$ws.Cells.Item(115,"U").Value = '28/10/2023 12:59'
$ws.Cells.Item(115,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/swinoujscie-unia-swarzedz/QD0sSekq/'

# Row 116
$ws.Cells.Item(116,"A").Value = 115
$ws.Cells.Item(116,"B").Value = 'poland'
$ws.Cells.Item(116,"C").Value = 'iii-liga-group-ii'
$ws.Cells.Item(116,"D").Value = '2023-2024'
$ws.Cells.Item(116,"E").Value = 45227.58333333334
$ws.Cells.Item(116,"F").Value = 'Elana Torun'
$ws.Cells.Item(116,"G").Value = 1
$ws.Cells.Item(116,"H").Value = 'Pogon Szczecin II'
$ws.Cells.Item(116,"I").Value = 0
$ws.Cells.Item(116,"J").Value = 2.24
$ws.Cells.Item(116,"K").Value = '27/10/2023 02:13'
$ws.Cells.Item(116,"L").Value = 2.15
$ws.Cells.Item(116,"M").Value = '28/10/2023 13:51'
$ws.Cells.Item(116,"N").Value = 3.26
$ws.Cells.Item(116,"O").Value = '27/10/2023 02:13'
$ws.Cells.Item(116,"P").Value = 3.37
$ws.Cells.Item(116,"Q").Value = '28/10/2023 13:51'
$ws.Cells.Item(116,"R").Value = 2.54
$ws.Cells.Item(116,"S").Value = '27/10/2023 02:13'
$ws.Cells.Item(116,"T").Value = 2.92
$ws.Cells.Item(116,"U").Value = '28/10/2023 13:51'
$ws.Cells.Item(116,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/elana-torun-pogon-szczecin/UsFBXHRL/'

# Row 117
$ws.Cells.Item(117,"A").Value = 116
$ws.Cells.Item(117,"B").Value = 'poland'
$ws.Cells.Item(117,"C").Value = 'iii-liga-group-ii'
$ws.Cells.Item(117,"D").Value = '2023-2024'
$ws.Cells.Item(117,"E").Value = 45227.58333333334
$ws.Cells.Item(117,"F").Value = 'Kleczew'
$ws.Cells.Item(117,"G").Value = 2
$ws.Cells.Item(117,"H").Value = 'Vineta W.'
$ws.Cells.Item(117,"I").Value = 1
$ws.Cells.Item(117,"J").Value = 2.1
$ws.Cells.Item(117,"K").Value = '27/10/2023 02:13'
$ws.Cells.Item(117,"L").Value = 2.3
$ws.Cells.Item(117,"M").Value = '28/10/2023 13:55'
$ws.Cells.Item(117,"N").Value = 3.36
$ws.Cells.Item(117,"O").Value = '27/10/2023 02:13'
$ws.Cells.Item(117,"P").Value = 3.72
$ws.Cells.Item(117,"Q").Value = '28/10/2023 13:55'
$ws.Cells.Item(117,"R").Value = 2.68
$ws.Cells.Item(117,"S").Value = '27/10/2023 02:13'
$ws.Cells.Item(117,"T").Value = 2.49
$ws.Cells.Item(117,"U").Value = '28/10/2023 13:55'
$ws.Cells.Item(117,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/kleczew-vineta-wolin/pfu8Ddzj/'

# Row 118
$ws.Cells.Item(118,"A").Value = 117
$ws.Cells.Item(118,"B").Value = 'poland'
$ws.Cells.Item(118,"C").Value = 'iii-liga-group-ii'
$ws.Cells.Item(118,"D").Value = '2023-2024'
$ws.Cells.Item(118,"E").Value = 45227.58333333334
$ws.Cells.Item(118,"F").Value = 'Notec Czarnkow'
$ws.Cells.Item(118,"G").Value = 1
$ws.Cells.Item(118,"H").Value = 'Nowe Skalmierzyce'
$ws.Cells.Item(118,"I").Value = 1
$ws.Cells.Item(118,"J").Value = 1.92
$ws.Cells.Item(118,"K").Value = '27/10/2023 02:13'
$ws.Cells.Item(118,"L").Value = 1.41
$ws.Cells.Item(118,"M").Value = '28/10/2023 13:58'
$ws.Cells.Item(118,"N").Value = 3.55
$ws.Cells.Item(118,"O").Value = '27/10/2023 02:13'
$ws.Cells.Item(118,"P").Value = 4.77
$ws.Cells.Item(118,"Q").Value = '28/10/2023 13:58'
$ws.Cells.Item(118,"R").Value = 2.9
$ws.Cells.Item(118,"S").Value = '27/10/2023 02:13'
$ws.Cells.Item(118,"T").Value = 5.35
$ws.Cells.Item(118,"U").Value = '28/10/2023 13:58'
$ws.Cells.Item(118,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/notec-czarnkow-nowe-skalmierzyce/YFL2Zwc9/'

# Row 119
$ws.Cells.Item(119,"A").Value = 118
$ws.Cells.Item(119,"B").Value = 'poland'
$ws.Cells.Item(119,"C").Value = 'iii-liga-group-ii'
$ws.Cells.Item(119,"D").Value = '2023-2024'
$ws.Cells.Item(119,"E").Value = 45227.60416666666
$ws.Cells.Item(119,"F").Value = 'Luzino'
$ws.Cells.Item(119,"G").Value = 3
$ws.Cells.Item(119,"H").Value = 'Starogard Gdanski'
$ws.Cells.Item(119,"I").Value = 0
$ws.Cells.Item(119,"J").Value = 2.36
$ws.Cells.Item(119,"K").Value = '27/10/2023 02:42'
$ws.Cells.Item(119,"L").Value = 2.31
$ws.Cells.Item(119,"M").Value = '28/10/2023 14:20'
$ws.Cells.Item(119,"N").Value = 3.33
$ws.Cells.Item(119,"O").Value = '27/10/2023 02:42'
$ws.Cells.Item(119,"P").Value = 3.44
$ws.Cells.Item(119,"Q").Value = '28/10/2023 14:20'
$ws.Cells.Item(119,"R").Value = 2.36
$ws.Cells.Item(119,"S").Value = '27/10/2023 02:42'
$ws.Cells.Item(119,"T").Value = 2.62
$ws.Cells.Item(119,"U").Value = '28/10/2023 14:20'
$ws.Cells.Item(119,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/luzino-starogard-gdanski/C6K6YcCF/'

# Row 120
$ws.Cells.Item(120,"A").Value = 119
$ws.Cells.Item(120,"B").Value = 'poland'
$ws.Cells.Item(120,"C").Value = 'iii-liga-group-ii'
$ws.Cells.Item(120,"D").Value = '2023-2024'
$ws.Cells.Item(120,"E").Value = 45227.625
$ws.Cells.Item(120,"F").Value = 'Sroda'
$ws.Cells.Item(120,"G").Value = 0
$ws.Cells.Item(120,"H").Value = 'Zawisza'
$ws.Cells.Item(120,"I").Value = 2
$ws.Cells.Item(120,"J").Value = 2.72
$ws.Cells.Item(120,"K").Value = '27/10/2023 03:12'
$ws.Cells.Item(120,"L").Value = 3.37
$ws.Cells.Item(120,"M").Value = '28/10/2023 14:56'
$ws.Cells.Item(120,"N").Value = 3.45
$ws.Cells.Item(120,"O").Value = '27/10/2023 03:12'
$ws.Cells.Item(120,"P").Value = 3.75
$ws.Cells.Item(120,"Q").Value = '28/10/2023 14:56'
$ws.Cells.Item(120,"R").Value = 2.04
$ws.Cells.Item(120,"S").Value = '27/10/2023 03:12'
$ws.Cells.Item(120,"T").Value = 1.84
$ws.Cells.Item(120,"U").Value = '28/10/2023 14:56'
$ws.Cells.Item(120,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/polonia-sroda-wielkopol-zawisza/OpvCCGkd/'

# Row 121
$ws.Cells.Item(121,"A").Value = 120
$ws.Cells.Item(121,"B").Value = 'poland'
$ws.Cells.Item(121,"C").Value = 'iii-liga-group-ii'
$ws.Cells.Item(121,"D").Value = '2023-2024'
$ws.Cells.Item(121,"E").Value = 45230.5625
$ws.Cells.Item(121,"F").Value = 'Cartusia Kartuzy'
$ws.Cells.Item(121,"G").Value = 2
$ws.Cells.Item(121,"H").Value = 'Elana Torun'
$ws.Cells.Item(121,"I").Value = 2
$ws.Cells.Item(121,"J").Value = 2.26
$ws.Cells.Item(121,"K").Value = '20/10/2023 01:12'
$ws.Cells.Item(121,"L").Value = 2.37
$ws.Cells.Item(121,"M").Value = '31/10/2023 13:08'
$ws.Cells.Item(121,"N").Value = 3.15
$ws.Cells.Item(121,"O").Value = '20/10/2023 01:12'
$ws.Cells.Item(121,"P").Value = 3.05
$ws.Cells.Item(121,"Q").Value = '31/10/2023 13:09'
$ws.Cells.Item(121,"R").Value = 2.59
$ws.Cells.Item(121,"S").Value = '20/10/2023 01:12'
$ws.Cells.Item(121,"T").Value = 2.82
$ws.Cells.Item(121,"U").Value = '31/10/2023 13:08'
$ws.Cells.Item(121,"V").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/cartusia-kartuzy-elana-torun/2ghHNtbc/'
